# Update the cryptos worksheet with the latest scraped symbol/price data
# (GitHub Actions run, Thu Dec 15 15:55:27 UTC 2022).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text; force text formatting before writing so
# numeric-looking strings (leading/trailing zeros, tiny decimals) are kept verbatim.
function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue "D2" "263.54"
Set-TextValue "D3" "22.71"
Set-TextValue "D4" "6.216"
Set-TextValue "D5" "0.06082"
Set-TextValue "D6" "3.517"
Set-TextValue "D7" "6.722"
Set-TextValue "D9" "0.8011"
Set-TextValue "D10" "0.1590"
Set-TextValue "D11" "0.08145"
Set-TextValue "D12" "0.03362"
Set-TextValue "D13" "0.03159"
Set-TextValue "D14" "0.09258"
Set-TextValue "D15" "3.921"
Set-TextValue "D16" "0.001684"
Set-TextValue "D17" "0.04825"
Set-TextValue "D18" "0.0006243"
Set-TextValue "D19" "0.006197"
Set-TextValue "D20" "0.005997"
$ws.Range("E20").Value = "19HotbitTokenHTB"
Set-TextValue "D21" "0.001099"
Set-TextValue "D24" "2.262"
Set-TextValue "D25" "0.3389"
Set-TextValue "D26" "0.1269"
Set-TextValue "D27" "0.0002687"
Set-TextValue "D40" "0.04636"
Set-TextValue "D41" "0.007275"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003908"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1118"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "D44" "0.01053"
Set-TextValue "D45" "0.00006111"
Set-TextValue "D46" "0.00000000752"
Set-TextValue "D47" "0.7516"
Set-TextValue "D48" "0.1593"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
